$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.010202050209045
$ws.Range("B1").Value = 2.12084698677063
$ws.Range("C1").Value = 6.246166706085205
$ws.Range("D1").Value = 1.451396703720093
$ws.Range("E1").Value = 1.338977336883545
